# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 508
$ws1.Range("F4").Value = 186
$ws1.Range("F7").Value = 90
$ws1.Range("F9").Value = 36
$ws1.Range("F10").Value = 6486
$ws1.Range("F11").Value = 219
$ws1.Range("F12").Value = 350
$ws1.Range("F13").Value = 2641
$ws1.Range("F14").Value = 158
$ws1.Range("F15").Value = 270
$ws1.Range("F17").Value = 504

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 508
$ws4.Range("F6").Value = 186
$ws4.Range("F9").Value = 90
$ws4.Range("F11").Value = 36
$ws4.Range("F13").Value = 6486
$ws4.Range("F15").Value = 219
$ws4.Range("F16").Value = 350
$ws4.Range("F17").Value = 2641
$ws4.Range("F18").Value = 158
$ws4.Range("F19").Value = 270
$ws4.Range("F21").Value = 504

$wb.Save()
